# Update countries & provincias Spain
# This script updates a handful of country statistics rows. Because the
# underlying data is kept sorted (descending) by "Casos totales" (column B),
# a couple of countries that grew past their neighbour's total now swap
# places with the row that used to sit above them. We therefore also
# relabel column A for the rows that shifted, while leaving their numeric
# data untouched (only the two "promoted" countries get fresh numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Israel (row 26): updated figures, no reordering needed -----------------
$ws.Range("B26").Value = 15466
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 6796
$ws.Range("E26").Value = 8468
$ws.Range("F26").Value = 129
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 202

# --- Ucrania overtakes Indonesia (rows 39-40) --------------------------------
# Row 39 used to be Indonesia; Ucrania moves up here with new numbers.
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 9009
$ws.Range("C39").Value = 392
$ws.Range("D39").Value = 864
$ws.Range("E39").Value = 7925
$ws.Range("F39").Value = 110
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 220

# Row 40 used to be Ucrania; Indonesia moves down here keeping its old numbers.
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 8882
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 1107
$ws.Range("E40").Value = 7032
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 743

# --- Kazajistan (row 62): updated figures ------------------------------------
$ws.Range("B62").Value = 2780
$ws.Range("C62").Value = 63
$ws.Range("E62").Value = 2073

# --- Georgia (row 108): updated figures --------------------------------------
$ws.Range("B108").Value = 496
$ws.Range("C108").Value = 10
$ws.Range("E108").Value = 341

# --- El Salvador overtakes Montenegro / Isla de Man / Tanzania (rows 122-125)
# Row 122 used to be Montenegro; El Salvador moves up here with new numbers.
$ws.Range("A122").Value = "El Salvador"
$ws.Range("B122").Value = 323
$ws.Range("C122").Value = 25
$ws.Range("D122").Value = 89
$ws.Range("E122").Value = 226
$ws.Range("F122").Value = 4
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 8

# Row 123 used to be Isla de Man; Montenegro moves down here, numbers unchanged.
$ws.Range("A123").Value = "Montenegro"
$ws.Range("B123").Value = 321
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 153
$ws.Range("E123").Value = 161
$ws.Range("F123").Value = 7
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 7

# Row 124 used to be Tanzania; Isla de Man moves down here, numbers unchanged.
$ws.Range("A124").Value = "Isla de Man"
$ws.Range("B124").Value = 308
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 247
$ws.Range("E124").Value = 43
$ws.Range("F124").Value = 22
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 18

# Row 125 used to be El Salvador; Tanzania moves down here, numbers unchanged.
$ws.Range("A125").Value = "Tanzania"
$ws.Range("B125").Value = 299
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 48
$ws.Range("E125").Value = 241
$ws.Range("F125").Value = 7
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 10

# --- Polinesia Francesa (row 165): updated figures ---------------------------
$ws.Range("D165").Value = 43
$ws.Range("E165").Value = 14

$wb.Save()
